# Mise à jour de l'application
# Adds a new "training day" column (BO) to the attendance sheet, mirroring
# the existing BN column's formatting, then updates the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcCol = 66   # BN
$newCol = 67   # BO

# Row 1: new date header (21/10/2025 -> serial 45951), formatted like BN1.
# (value first, formatting copy second -- doing the PasteSpecial before the
# value write leaves dependent COUNTA/COUNTIF formulas with a stale cache)
$ws.Cells.Item(1, $newCol).Value2 = 45951
$ws.Cells.Item(1, $srcCol).Copy()
$ws.Cells.Item(1, $newCol).PasteSpecial(-4122)

# Attendance letter per player row, copying BN's style onto BO first so the
# new cell matches the rest of the table, then writing the value.
$rowValues = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "RH"
    18 = "P"
    19 = "B"
    20 = "P"
    21 = "B"
    22 = "P"
    23 = "B"
    24 = "B"
    25 = "P"
    26 = "P"
    27 = "P"
    28 = "P"
    29 = "P"
}

foreach ($row in 2..29) {
    if ($rowValues.ContainsKey($row)) {
        $ws.Cells.Item($row, $newCol).Value2 = $rowValues[$row]
        $ws.Cells.Item($row, $srcCol).Copy()
        $ws.Cells.Item($row, $newCol).PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = $false

# Update the frozen-pane scroll position and the active selection to match
# the latest columns now that BO has been added.
$ws.Range("BM1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("BQ26").Select()

$wb.Save()
